$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge first_name (B) and last_name (C) columns into a single "name" column (B),
# shifting skills/experience/contact left by one column.

# Update data first (row 2) before headers, order doesn't really matter here,
# but do data first then shift headers, then delete the now-empty last_name column.

$ws.Range("B1").Value = "name"
$ws.Range("B2").Value = "Omar Rodriguez-Lopez"

# Remove the now redundant last_name column (old column C), which shifts
# skills/experience/contact left into C/D/E.
$ws.Columns("C").Delete()

# Set column widths to match target formatting (values chosen so that,
# after Excel's internal pixel-rounding of ColumnWidth, the persisted
# width attribute lands as close as possible to the target).
$ws.Columns("B").ColumnWidth = 33.83
$ws.Columns("C").ColumnWidth = 61.0

$ws.Range("C14").Select()
